# The deck's applied design ("Integral") is being reverted back to the
# default Office colour theme (the notes master keeps the original
# "Office Theme" colours while the slide master/presentation theme was
# using "Integral" — this edit swaps the two, bringing the presentation
# back to plain Office colours).
#
# PowerPoint's theme colours are changed through the ThemeColorScheme /
# ThemeColor object model (Slide.ThemeColorScheme), which maps 1:1 onto
# the <a:clrScheme> children in order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
